$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. CORE COMPETENCIES: collapse the three detailed bullet paragraphs into a
#    single short summary line ("Product Marketing Core (bullet) Research &
#    Analytics (bullet) Communication & Technology").
# ---------------------------------------------------------------------------
$bullet = [char]0x2022

$coreP = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Contains("Product Marketing Core:")) {
        $coreP = $d.Paragraphs($i)
        break
    }
}

$researchP = $coreP.Next()
$commsP = $researchP.Next()

# Remove the second and third paragraphs (and their paragraph marks) entirely.
$removeStart = $researchP.Range.Start
$removeEnd = $commsP.Range.End
$d.Range($removeStart, $removeEnd).Delete()

# Replace the remaining (first) paragraph's text with the condensed summary.
$coreP.Range.Text = "Product Marketing Core " + $bullet + " Research & Analytics " + $bullet + " Communication & Technology"

# ---------------------------------------------------------------------------
# 2. Add a new "TECHNICAL SKILLS" section (heading + 3 detail paragraphs,
#    using the text that used to live under CORE COMPETENCIES) right after
#    the last bullet of Key Achievements / Data-Driven Marketing, before the
#    closing "For a more detailed..." paragraph.
# ---------------------------------------------------------------------------
$anchorP = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Contains("Managed national polling team of five data analysts")) {
        $anchorP = $d.Paragraphs($i)
        break
    }
}

# Insert four empty paragraphs directly after the anchor; each call inserts
# immediately after $anchorP, so the last-inserted ends up closest to it and
# the first-inserted (done first) ends up furthest -- i.e. the net order
# below matches the order the statements were issued.
$anchorP.Range.InsertParagraphAfter()
$anchorP.Range.InsertParagraphAfter()
$anchorP.Range.InsertParagraphAfter()
$anchorP.Range.InsertParagraphAfter()

$headingP = $anchorP.Next()
$coreSkillsP = $headingP.Next()
$researchSkillsP = $coreSkillsP.Next()
$commsSkillsP = $researchSkillsP.Next()

$headingP.Range.Text = "TECHNICAL SKILLS"
$headingP.Style = "Heading 2"

$coreSkillsP.Range.Text = "PRODUCT MARKETING CORE Market Intelligence & Competitive Analysis; Product Positioning & Messaging Development; Go-to-Market Strategy & Product Launch Management; Customer Segmentation & Buyer Persona Development"

$researchSkillsP.Range.Text = "RESEARCH & ANALYTICS Survey Methodology & Customer Insights; Market Research Design & Implementation; Competitive Intelligence & SWOT Analysis; A/B Testing & Conversion Optimization"

$commsSkillsP.Range.Text = "COMMUNICATION & TECHNOLOGY Strategic Messaging & Narrative Development; Technical Concept Translation for Business Audiences; Data Visualization & Reporting (Tableau, PowerBI, d3.js); Client Relationship Management & Business Development"
